$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-03 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-04 Wednesday", 2)

$d.Content.Find.Execute("76×46=3496", $true, $false, $false, $false, $false, $true, 1, $false, "37×73=2701", 2)
$d.Content.Find.Execute("53×67=3551", $true, $false, $false, $false, $false, $true, 1, $false, "74×44=3256", 2)
$d.Content.Find.Execute("39×94=3666", $true, $false, $false, $false, $false, $true, 1, $false, "60×28=1680", 2)
$d.Content.Find.Execute("50×64=3200", $true, $false, $false, $false, $false, $true, 1, $false, "26×13=338", 2)
$d.Content.Find.Execute("51×51=2601", $true, $false, $false, $false, $false, $true, 1, $false, "32×77=2464", 2)

$d.Content.Find.Execute("76×67=5092", $true, $false, $false, $false, $false, $true, 1, $false, "50×72=3600", 2)
$d.Content.Find.Execute("50×67=3350", $true, $false, $false, $false, $false, $true, 1, $false, "41×86=3526", 2)
$d.Content.Find.Execute("89×69=6141", $true, $false, $false, $false, $false, $true, 1, $false, "53×25=1325", 2)
$d.Content.Find.Execute("70×70=4900", $true, $false, $false, $false, $false, $true, 1, $false, "98×79=7742", 2)
$d.Content.Find.Execute("64×12=768", $true, $false, $false, $false, $false, $true, 1, $false, "19×80=1520", 2)

$d.Content.Find.Execute("26×90=2340", $true, $false, $false, $false, $false, $true, 1, $false, "96×81=7776", 2)
$d.Content.Find.Execute("26×45=1170", $true, $false, $false, $false, $false, $true, 1, $false, "20×37=740", 2)
$d.Content.Find.Execute("54×68=3672", $true, $false, $false, $false, $false, $true, 1, $false, "76×42=3192", 2)
$d.Content.Find.Execute("68×31=2108", $true, $false, $false, $false, $false, $true, 1, $false, "67×25=1675", 2)
$d.Content.Find.Execute("11×95=1045", $true, $false, $false, $false, $false, $true, 1, $false, "88×17=1496", 2)

$d.Content.Find.Execute("36×79=2844", $true, $false, $false, $false, $false, $true, 1, $false, "72×54=3888", 2)
$d.Content.Find.Execute("63×72=4536", $true, $false, $false, $false, $false, $true, 1, $false, "95×63=5985", 2)
$d.Content.Find.Execute("28×17=476", $true, $false, $false, $false, $false, $true, 1, $false, "14×59=826", 2)
$d.Content.Find.Execute("81×40=3240", $true, $false, $false, $false, $false, $true, 1, $false, "28×65=1820", 2)
$d.Content.Find.Execute("12×93=1116", $true, $false, $false, $false, $false, $true, 1, $false, "81×92=7452", 2)

$d.Content.Find.Execute("37×14=518", $true, $false, $false, $false, $false, $true, 1, $false, "36×50=1800", 2)
$d.Content.Find.Execute("76×69=5244", $true, $false, $false, $false, $false, $true, 1, $false, "26×67=1742", 2)
$d.Content.Find.Execute("70×98=6860", $true, $false, $false, $false, $false, $true, 1, $false, "93×98=9114", 2)
$d.Content.Find.Execute("34×40=1360", $true, $false, $false, $false, $false, $true, 1, $false, "13×83=1079", 2)
$d.Content.Find.Execute("99×30=2970", $true, $false, $false, $false, $false, $true, 1, $false, "69×13=897", 2)
